# Generate Report for Handoff
#
# The localization status report is regenerated: the file
# "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md" has progressed from
# "In Translation" to "Ready for handoff" (with a fresh handoff
# timestamp), while "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md" stays
# "In Translation". As a result the two data rows on every sheet swap
# places (ecbc5422 now sorts first, 66fc1641 second).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value2 = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$ov.Range("B2").Value2 = "In Translation"
$ov.Range("C2").Value2 = "In Translation"
$ov.Range("D2").Value2 = "2016-03-21 08:16:45"

$ov.Range("A3").Value2 = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$ov.Range("B3").Value2 = "Ready for handoff"
$ov.Range("C3").Value2 = "Ready for handoff"
$ov.Range("D3").Value2 = "2016-03-21 08:17:33"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", [Type]::Missing, [Type]::Missing, "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", [Type]::Missing, [Type]::Missing, "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value2 = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$zh.Range("B2").Value2 = ".md"
$zh.Range("C2").Value2 = "In Translation"
$zh.Range("D2").Value2 = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf"
$zh.Range("E2").Value2 = "2016-03-21 08:16:41"
$zh.Range("H2").Value2 = "0001-01-01 00:00:00"
$zh.Range("J2").Value2 = "Include"

$zh.Range("A3").Value2 = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$zh.Range("B3").Value2 = ".md"
$zh.Range("C3").Value2 = "Ready for handoff"
$zh.Range("D3").Value2 = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf"
$zh.Range("E3").Value2 = "2016-03-21 08:17:29"
$zh.Range("H3").Value2 = "0001-01-01 00:00:00"
$zh.Range("J3").Value2 = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", [Type]::Missing, [Type]::Missing, "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25e1da1c23593d68aa60eac33d03fe5ec62c421e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", [Type]::Missing, [Type]::Missing, "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25e1da1c23593d68aa60eac33d03fe5ec62c421e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value2 = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md"
$de.Range("B2").Value2 = ".md"
$de.Range("C2").Value2 = "In Translation"
$de.Range("D2").Value2 = "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf"
$de.Range("E2").Value2 = "2016-03-21 08:16:45"
$de.Range("H2").Value2 = "0001-01-01 00:00:00"
$de.Range("J2").Value2 = "Include"

$de.Range("A3").Value2 = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md"
$de.Range("B3").Value2 = ".md"
$de.Range("C3").Value2 = "Ready for handoff"
$de.Range("D3").Value2 = "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf"
$de.Range("E3").Value2 = "2016-03-21 08:17:33"
$de.Range("H3").Value2 = "0001-01-01 00:00:00"
$de.Range("J3").Value2 = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md", [Type]::Missing, [Type]::Missing, "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daacd1925ee868f8937a97b32af7b49aedf83c28/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf", [Type]::Missing, [Type]::Missing, "ecbc5422-8f00-4f03-8bcb-e412fcafa55a.cc9ff078b38cfa5e35f4b05cc10d6e4feb82cd40.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce1661b5762c072c659b9a36cd06837c06851ce1/e2e/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md", [Type]::Missing, [Type]::Missing, "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daacd1925ee868f8937a97b32af7b49aedf83c28/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf", [Type]::Missing, [Type]::Missing, "66fc1641-57ea-4ee9-8aa3-e5a139707f2f.2d10009bacfec9e590dd5412141c0a48f18e9a5a.de-de.xlf") | Out-Null

$wb.Save()
